# ActivityLogSheetWk5.xlsx - fill in the Week 5 activity log entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (row 2) ---
# B2:E2 is merged - write the student's name into the anchor cell
$ws.Range("B2").Value = "Jesse Hare"
# F2 used to hold the static label "Week"; it now shows "Week 5" directly
$ws.Range("F2").Value = "Week 5"

# --- Activity rows (4-8); A:B is merged per row, pre-existing styles/number
#     formats on C:G already match the target (date / time / text), so we
#     only need to drop values into the template's blank cells. ---
$ws.Range("A4").Value = "Design Screen Layouts"
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = 43703
$ws.Range("E4").Value = 0.375
$ws.Range("F4").Value = 0.41666666666666669
$ws.Range("G4").Value = "Group"

$ws.Range("A5").Value = "Finalise framework choice"
$ws.Range("C5").Value = "G"
$ws.Range("D5").Value = 43704
$ws.Range("E5").Value = 0.41666666666666669
$ws.Range("F5").Value = 0.45833333333333331
$ws.Range("G5").Value = "Group"

$ws.Range("A6").Value = "Identify needed libraries"
$ws.Range("C6").Value = "G"
$ws.Range("D6").Value = 43705
$ws.Range("E6").Value = 0.54166666666666663
$ws.Range("F6").Value = 0.58333333333333337
$ws.Range("G6").Value = "Group"

$ws.Range("A7").Value = "Identify needed libraries"
$ws.Range("C7").Value = "G"
$ws.Range("D7").Value = 43706
$ws.Range("E7").Value = 0.35416666666666669
$ws.Range("F7").Value = 0.39583333333333331
$ws.Range("G7").Value = "Group"

$ws.Range("A8").Value = "Identify needed libraries"
$ws.Range("C8").Value = "G"
$ws.Range("D8").Value = 43707
$ws.Range("E8").Value = 0.38541666666666669
$ws.Range("F8").Value = 0.42708333333333331
$ws.Range("G8").Value = "Group"

# --- Totals row: G17 used to be =SUM(G4:G16); the group-hours column now
#     holds text labels instead of numbers, so the author replaced the
#     formula with a hand-entered total. ---
$ws.Range("G17").Value = 20

# --- Column width tweaks (values chosen so the stored/quantized width in
#     the saved file lands on the same figures as the target file) ---
$ws.Columns.Item(2).ColumnWidth = 13.15   # -> stored width 14
$ws.Columns.Item(4).ColumnWidth = 12.15   # -> stored width 13
$ws.Columns.Item(5).ColumnWidth = 12.75   # -> stored width ~13.7 (best-fit)
$ws.Columns.Item(6).ColumnWidth = 12.75   # -> stored width ~13.7 (best-fit)

# --- Selection cosmetics, matches the saved view in the target file ---
$ws.Range("A5:B5").Select()
